$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 1772.6364
$ws.Cells.Item(17, 10).Value = 1772.6364
$ws.Cells.Item(17, 12).Value = 5317.9092
$ws.Cells.Item(17, 14).Value = -5653.9092

$ws.Cells.Item(40, 8).Value = 4036.5
$ws.Cells.Item(40, 9).Value = 2987.5
$ws.Cells.Item(40, 10).Value = 4561
$ws.Cells.Item(40, 11).Value = 2987.5
$ws.Cells.Item(40, 12).Value = 4561
$ws.Cells.Item(40, 13).Value = -2812.5
$ws.Cells.Item(40, 14).Value = -4911

$ws.Cells.Item(55, 8).Value = 2263.5386
$ws.Cells.Item(55, 9).Value = 332.5
$ws.Cells.Item(55, 11).Value = 332.5
$ws.Cells.Item(55, 13).Value = -118.5

$ws.Cells.Item(64, 8).Value = 4856.4346
$ws.Cells.Item(64, 9).Value = 4339.8
$ws.Cells.Item(64, 10).Value = 4999.9443
$ws.Cells.Item(64, 11).Value = 4339.8
$ws.Cells.Item(64, 12).Value = 4999.9443
$ws.Cells.Item(64, 13).Value = -4091.8
$ws.Cells.Item(64, 14).Value = -5495.9443

$ws.Cells.Item(67, 8).Value = 4856.4346
$ws.Cells.Item(67, 9).Value = 4339.8
$ws.Cells.Item(67, 10).Value = 4999.9443
$ws.Cells.Item(67, 11).Value = 4339.8
$ws.Cells.Item(67, 12).Value = 4999.9443
$ws.Cells.Item(67, 13).Value = -3481.8
$ws.Cells.Item(67, 14).Value = -6715.9443

$ws.Cells.Item(132, 8).Value = 2200.4119
$ws.Cells.Item(132, 9).Value = 1902.375
$ws.Cells.Item(132, 11).Value = 5707.125
$ws.Cells.Item(132, 13).Value = -3177.125

$ws.Cells.Item(138, 8).Value = 2690.875
$ws.Cells.Item(138, 9).Value = 1922.4286
$ws.Cells.Item(138, 10).Value = 2785.2456
$ws.Cells.Item(138, 11).Value = 5767.2858
$ws.Cells.Item(138, 12).Value = 8355.736800000001
$ws.Cells.Item(138, 13).Value = -627.2857999999997
$ws.Cells.Item(138, 14).Value = -18635.7368

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 13192927
$ws.Cells.Item(32, 9).Value = 16705507
$ws.Cells.Item(32, 11).Value = 16705507
$ws.Cells.Item(32, 13).Value = -16705220

$ws.Cells.Item(132, 8).Value = 3471.52
$ws.Cells.Item(132, 9).Value = 2091.0322
$ws.Cells.Item(132, 11).Value = 6273.096600000001
$ws.Cells.Item(132, 13).Value = -3743.096600000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(123, 8).Value = 81041.60000000001
$ws.Cells.Item(123, 10).Value = 81041.60000000001
$ws.Cells.Item(123, 12).Value = 81041.60000000001
$ws.Cells.Item(123, 14).Value = -90841.60000000001

$ws.Cells.Item(134, 8).Value = 278894.75
$ws.Cells.Item(134, 9).Value = 1168.6875
$ws.Cells.Item(134, 11).Value = 3506.0625
$ws.Cells.Item(134, 13).Value = -971.0625

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 2756.353
$ws.Cells.Item(58, 9).Value = 931.3333
$ws.Cells.Item(58, 11).Value = 931.3333
$ws.Cells.Item(58, 13).Value = -728.3333

$ws.Cells.Item(107, 8).Value = 1981.826
$ws.Cells.Item(107, 9).Value = 1024.6364
$ws.Cells.Item(107, 11).Value = 1024.6364
$ws.Cells.Item(107, 13).Value = 895.3635999999999

$ws.Cells.Item(134, 8).Value = 3160.7715
$ws.Cells.Item(134, 9).Value = 1885.2
$ws.Cells.Item(134, 11).Value = 5655.6
$ws.Cells.Item(134, 13).Value = -3120.6

$ws.Cells.Item(136, 8).Value = 2756.353
$ws.Cells.Item(136, 9).Value = 931.3333
$ws.Cells.Item(136, 11).Value = 2793.9999
$ws.Cells.Item(136, 13).Value = -243.9998999999998

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(29, 8).Value = 116.57143
$ws.Cells.Item(29, 9).Value = 57.25
$ws.Cells.Item(29, 10).Value = 195.66667
$ws.Cells.Item(29, 11).Value = 171.75
$ws.Cells.Item(29, 12).Value = 587.00001
$ws.Cells.Item(29, 13).Value = 105.25
$ws.Cells.Item(29, 14).Value = -1141.00001

$ws.Cells.Item(38, 8).Value = 291.1111
$ws.Cells.Item(38, 9).Value = 696.3333
$ws.Cells.Item(38, 10).Value = 88.5
$ws.Cells.Item(38, 11).Value = 2088.9999
$ws.Cells.Item(38, 12).Value = 265.5
$ws.Cells.Item(38, 13).Value = -1741.9999
$ws.Cells.Item(38, 14).Value = -959.5

$ws.Cells.Item(46, 8).Value = 1647.8334
$ws.Cells.Item(46, 9).Value = 971.75
$ws.Cells.Item(46, 11).Value = 2915.25
$ws.Cells.Item(46, 13).Value = -2824.25

$ws.Cells.Item(107, 8).Value = 581.7143
$ws.Cells.Item(107, 9).Value = 427.4762
$ws.Cells.Item(107, 10).Value = 1044.4286
$ws.Cells.Item(107, 11).Value = 1282.4286
$ws.Cells.Item(107, 12).Value = 3133.2858
$ws.Cells.Item(107, 13).Value = 637.5714
$ws.Cells.Item(107, 14).Value = -6973.2858

$ws.Cells.Item(118, 8).Value = 6000
$ws.Cells.Item(118, 9).Value = 6000
$ws.Cells.Item(118, 10).Value = 0
$ws.Cells.Item(118, 11).Value = 18000
$ws.Cells.Item(118, 12).Value = 0
$ws.Cells.Item(118, 13).Value = -16757
$ws.Cells.Item(118, 14).ClearContents()

$ws.Cells.Item(131, 8).Value = 5103.607
$ws.Cells.Item(131, 9).Value = 3952.6667
$ws.Cells.Item(131, 11).Value = 11858.0001
$ws.Cells.Item(131, 13).Value = -6818.000100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 7647.6665
$ws.Cells.Item(70, 9).Value = 6996.75
$ws.Cells.Item(70, 11).Value = 6996.75
$ws.Cells.Item(70, 13).Value = -6726.75

$ws.Cells.Item(73, 8).Value = 7647.6665
$ws.Cells.Item(73, 9).Value = 6996.75
$ws.Cells.Item(73, 11).Value = 6996.75
$ws.Cells.Item(73, 13).Value = -6060.75

$ws.Cells.Item(93, 8).Value = 59989.332
$ws.Cells.Item(93, 10).Value = 59989.332
$ws.Cells.Item(93, 12).Value = 59989.332
$ws.Cells.Item(93, 14).Value = -63733.332

$ws.Cells.Item(97, 8).Value = 1311.1765
$ws.Cells.Item(97, 9).Value = 1845.7778
$ws.Cells.Item(97, 11).Value = 1845.7778
$ws.Cells.Item(97, 13).Value = -1349.7778

$ws.Cells.Item(109, 8).Value = 45241.5
$ws.Cells.Item(109, 10).Value = 45241.5
$ws.Cells.Item(109, 12).Value = 45241.5
$ws.Cells.Item(109, 14).Value = -47321.5

$ws.Cells.Item(122, 8).Value = 2015.8462
$ws.Cells.Item(122, 9).Value = 1675.125
$ws.Cells.Item(122, 11).Value = 5025.375
$ws.Cells.Item(122, 13).Value = -2575.375

$ws.Cells.Item(132, 8).Value = 21281676
$ws.Cells.Item(132, 9).Value = 27780766
$ws.Cells.Item(132, 10).Value = 11922.272
$ws.Cells.Item(132, 11).Value = 83342298
$ws.Cells.Item(132, 12).Value = 35766.81600000001
$ws.Cells.Item(132, 13).Value = -83339768
$ws.Cells.Item(132, 14).Value = -40826.81600000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 118781.11
$ws.Cells.Item(7, 9).Value = 6003.6665
$ws.Cells.Item(7, 10).Value = 175169.83
$ws.Cells.Item(7, 11).Value = 6003.6665
$ws.Cells.Item(7, 12).Value = 175169.83
$ws.Cells.Item(7, 13).Value = -5891.6665
$ws.Cells.Item(7, 14).Value = -175393.83

$ws.Cells.Item(16, 8).Value = 495.4375
$ws.Cells.Item(16, 9).Value = 495.4375
$ws.Cells.Item(16, 11).Value = 495.4375
$ws.Cells.Item(16, 13).Value = -325.4375

$ws.Cells.Item(61, 8).Value = 3590.3845
$ws.Cells.Item(61, 9).Value = 3556.25
$ws.Cells.Item(61, 10).Value = 4000
$ws.Cells.Item(61, 11).Value = 3556.25
$ws.Cells.Item(61, 12).Value = 4000
$ws.Cells.Item(61, 13).Value = -3354.25
$ws.Cells.Item(61, 14).Value = -4404

$ws.Cells.Item(110, 8).Value = 57000
$ws.Cells.Item(110, 10).Value = 57000
$ws.Cells.Item(110, 12).Value = 57000
$ws.Cells.Item(110, 14).Value = -65180

$ws.Cells.Item(113, 8).Value = 3590.3845
$ws.Cells.Item(113, 9).Value = 3556.25
$ws.Cells.Item(113, 10).Value = 4000
$ws.Cells.Item(113, 11).Value = 3556.25
$ws.Cells.Item(113, 12).Value = 4000
$ws.Cells.Item(113, 13).Value = -1386.25
$ws.Cells.Item(113, 14).Value = -8340

$ws.Cells.Item(126, 8).Value = 118781.11
$ws.Cells.Item(126, 9).Value = 6003.6665
$ws.Cells.Item(126, 10).Value = 175169.83
$ws.Cells.Item(126, 11).Value = 18010.9995
$ws.Cells.Item(126, 12).Value = 525509.49
$ws.Cells.Item(126, 13).Value = -15540.9995
$ws.Cells.Item(126, 14).Value = -530449.49

$ws.Cells.Item(132, 8).Value = 487283.75
$ws.Cells.Item(132, 9).Value = 13704.5625
$ws.Cells.Item(132, 11).Value = 41113.6875
$ws.Cells.Item(132, 13).Value = -38583.6875

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(64, 8).Value = 65000
$ws.Cells.Item(64, 10).Value = 65000
$ws.Cells.Item(64, 12).Value = 65000
$ws.Cells.Item(64, 14).Value = -65496

$ws.Cells.Item(67, 8).Value = 65000
$ws.Cells.Item(67, 10).Value = 65000
$ws.Cells.Item(67, 12).Value = 65000
$ws.Cells.Item(67, 14).Value = -66716

$ws.Cells.Item(96, 8).Value = 4170.6
$ws.Cells.Item(96, 9).Value = 4061.111
$ws.Cells.Item(96, 11).Value = 4061.111
$ws.Cells.Item(96, 13).Value = -2688.111

$ws.Cells.Item(113, 8).Value = 1439.2
$ws.Cells.Item(113, 9).Value = 749.5
$ws.Cells.Item(113, 10).Value = 1899
$ws.Cells.Item(113, 11).Value = 2248.5
$ws.Cells.Item(113, 12).Value = 5697
$ws.Cells.Item(113, 13).Value = -78.5
$ws.Cells.Item(113, 14).Value = -10037

$ws.Cells.Item(132, 8).Value = 1983.3846
$ws.Cells.Item(132, 9).Value = 1322.4762
$ws.Cells.Item(132, 11).Value = 3967.4286
$ws.Cells.Item(132, 13).Value = -1437.4286

$ws.Cells.Item(135, 8).Value = 90342.8
$ws.Cells.Item(135, 10).Value = 90342.8
$ws.Cells.Item(135, 12).Value = 90342.8
$ws.Cells.Item(135, 14).Value = -100482.8
